$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view metadata ---
$excel.Windows.Item(1).WindowState = -4143  # xlNormal (no-op placeholder, harmless)

# --- Row 11: new "Model 10" data row ---
# Copy number formatting/style from row 10 (A:D) for the text/number cells
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0

# Copy style from F2 (centered, bordered, non-bold) for E11:H11, then drop vertical centering
$ws.Range("F2").Copy()
$ws.Range("E11:H11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0
$ws.Range("E11:H11").VerticalAlignment = -4107  # xlVAlignBottom (renders as "no vertical override")

# Set the row 11 values
$ws.Range("A11").Value = "Model 10"
$ws.Range("B11").Value = "Resnet-C"
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = "Model10_history"
$ws.Range("E11").Value = 0.3841
$ws.Range("F11").Value = 0.8472
$ws.Range("G11").Value = 0.4577
$ws.Range("H11").Value = 0.8304

# Give I11 the same formatting as I10 first, then add the hyperlink, then
# restore the exact formatting afterwards (Hyperlinks.Add resets font/style).
$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Hyperlinks.Add($ws.Range("I11"), "https://github.com/Coachnmomof3/UCB_COVID_Prediction_Model/blob/Gabriel_Cuchacovich/GoogleColab_COVID_ML.ipynb")

$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# --- Sheet view selection ---
$ws.Range("G32").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1  # xlPortrait
